$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$ws.Cells.Item(4, 1).Value = "test@test.com"
$ws.Cells.Item(4, 2).Value = "127.0.0.1"
$ws.Cells.Item(4, 3).Value = "Mozilla/5.0 (Windows NT 10.0; Win64; x64) AppleWebKit/537.36 (KHTML, like Gecko) Chrome/139.0.0.0 Safari/537.36 Edg/139.0.0.0"
$ws.Cells.Item(4, 4).Value = "2025-08-22T09:59:40.556063+00:00"
